# Apply "New PO forecast model" update
$wb = $excel.ActiveWorkbook

# --- Sheet "Weekly Quantity": append 3 new weekly rows ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$weeklyNewRows = @(
    @(45662.99999999999, 60),
    @(45669.99999999999, 20),
    @(45683.99999999999, 20)
)
$r = 65
foreach ($row in $weeklyNewRows) {
    $ws1.Range("A" + $r).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# --- Sheet "Monthly Trend": append 1 new monthly row ---
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("A25").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(25, 1).Value = 45688.99999999999
$ws2.Cells.Item(25, 2).Value = 100

# --- Sheet "PO Forecast": new forecast model ---
$ws3 = $wb.Worksheets.Item("PO Forecast")

# Rows 2-64 keep the same week-ending dates; only the forecast value changes.
$forecastValues = @(0,0,94,174,219,269,377,547,429,117,0,222,439,521,315,166,101,679,682,410,107,130,1253,894,343,1141,1042,527,10,0,956,915,410,0,128,262,349,399,444,539,702,858,870,356,156,364,595,705,666,891,658,334,278,1126,945,850,540,601,1283,266,102,474,688)
for ($i = 0; $i -lt $forecastValues.Length; $i++) {
    $ws3.Cells.Item($i + 2, 2).Value = $forecastValues[$i]
}

# Rows 65-75: replaced forecast tail (new dates + values), extending the sheet.
$forecastTail = @(
    @(45662.99999999999, 134),
    @(45669.99999999999, 287),
    @(45683.99999999999, 524),
    @(45690.99999999999, 578),
    @(45697.99999999999, 621),
    @(45704.99999999999, 703),
    @(45711.99999999999, 857),
    @(45718.99999999999, 1024),
    @(45725.99999999999, 1071),
    @(45732.99999999999, 905),
    @(45739.99999999999, 598)
)
$r = 65
foreach ($row in $forecastTail) {
    $ws3.Range("A" + $r).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

